$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing values ---
# Valor Mora
$ws.Range("E11").Value = 170820
# Cant. Trabajadores
$ws.Range("C13").Value = 2

# --- Insert a new data row (17) for the new worker, copying the format
#     of the existing "middle" data row (16) so the table keeps its
#     boxed border look ---
$ws.Rows.Item(17).Insert()
$ws.Range("B16:J16").Copy()
$ws.Range("B17:J17").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("B17").Value = "CC"
$ws.Range("C17").Value = "1047384474"
$ws.Range("D17").Value = "JORGE ELIECER ALVAREZ SAMPAYO"
$ws.Range("E17").Value = "2508"
$ws.Range("F17").Value = 56940
$ws.Range("G17").Value = 1423500

# --- The pre-existing second record (now on row 18 after the insert)
#     moves its "Periodo Mora" forward to the new period ---
$ws.Range("E18").Value = "2508"

# --- Re-fit the data columns now that a new (longer) name and a new row
#     of values were added to the table ---
$ws.Columns.Item(2).ColumnWidth = 16.90625
$ws.Columns.Item(3).ColumnWidth = 10.81640625
$ws.Columns.Item(4).ColumnWidth = 29.54296875
$ws.Columns.Item(5).ColumnWidth = 12.7265625
$ws.Columns.Item(6).ColumnWidth = 9.453125
$ws.Columns.Item(7).ColumnWidth = 13.453125
$ws.Columns.Item(8).ColumnWidth = 17.90625
$ws.Columns.Item(9).ColumnWidth = 16.81640625
$ws.Columns.Item(10).ColumnWidth = 14.1796875

# --- Nudge the logo a little to the left to re-center it over the
#     (now narrower) header columns ---
$logo = $ws.Shapes.Item(1)
$logoWidth = 975600 / 12700
$logoHeight = 612000 / 12700
$logo.Left = $logo.Left - 13.5
$logo.Width = $logoWidth
$logo.Height = $logoHeight
